$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds 14 data rows (rows 2-15): line1..line6, extr1..extr8.
# The edit inserts two new rows "line7" and "line8" right after "line6" (i.e. at
# sheet rows 8 and 9), which pushes extr1..extr8 down by two rows (to rows 10-17).
# A handful of the extr* rows' C/D/E values also change as part of the same commit.

# Step 1: capture the current rows 8-15 (extr1..extr8: name + A index) before
# they get overwritten, so we can re-write them two rows further down.
$oldNames = @()
for ($r = 8; $r -le 15; $r++) {
    $oldNames += , ($ws.Cells.Item($r, 2).Value())
}

# Step 2: new C/D/E values (post-edit) for extr1..extr8, in order.
$newCDE = @(
    @(5, 12, $true),   # extr1
    @(5, 9,  $true),   # extr2
    @(10, 11, $true),  # extr3
    @(7, 8,  $false),  # extr4
    @(9, 11, $true),   # extr5
    @(7, 11, $true),   # extr6
    @(5, 7,  $true),   # extr7
    @(8, 5,  $true)    # extr8
)

# Step 3: write extr1..extr8 into their new home, rows 10-17.
for ($i = 0; $i -lt 8; $i++) {
    $destRow = 10 + $i
    $ws.Cells.Item($destRow, 1).Value = $destRow - 2
    $ws.Cells.Item($destRow, 2).Value = $oldNames[$i]
    $ws.Cells.Item($destRow, 3).Value = $newCDE[$i][0]
    $ws.Cells.Item($destRow, 4).Value = $newCDE[$i][1]
    $ws.Cells.Item($destRow, 5).Value = $newCDE[$i][2]
}

# Step 4: rows 16 and 17 are brand-new sheet rows; give their A cell the same
# formatting (border/bold/alignment) as the rest of the index column by
# copying the format (only) from an existing, already-styled cell.
$ws.Range("A2").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Step 5: write the two new "line7"/"line8" rows into rows 8-9.
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $true
